$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Snapshot original values for rows 2-36 (columns D, M, N, O, P, R, S)
$snapD = @{}
$snapM = @{}
$snapN = @{}
$snapO = @{}
$snapP = @{}
$snapR = @{}
$snapS = @{}

for ($i = 2; $i -le 36; $i++) {
    $snapD[$i] = $ws.Cells.Item($i, 4).Value2
    $snapM[$i] = $ws.Cells.Item($i, 13).Value2
    $snapN[$i] = $ws.Cells.Item($i, 14).Value2
    $snapO[$i] = $ws.Cells.Item($i, 15).Value2
    $snapP[$i] = $ws.Cells.Item($i, 16).Value2
    $snapR[$i] = $ws.Cells.Item($i, 18).Value2
    $snapS[$i] = $ws.Cells.Item($i, 19).Value2
}

# Mapping of target row -> source row (permutation derived from the diff)
$map = @{}
$map[2] = 28
$map[3] = 9
$map[4] = 25
$map[5] = 14
$map[6] = 31
$map[7] = 35
$map[8] = 8
$map[9] = 16
$map[10] = 30
$map[11] = 36
$map[12] = 12
$map[13] = 33
$map[14] = 17
$map[15] = 6
$map[16] = 11
$map[17] = 20
$map[18] = 2
$map[19] = 13
$map[20] = 5
$map[21] = 21
$map[22] = 27
$map[23] = 3
$map[24] = 15
$map[25] = 23
$map[26] = 34
$map[27] = 18
$map[28] = 22
$map[29] = 32
$map[30] = 26
$map[31] = 19
$map[32] = 24
$map[33] = 10
$map[34] = 7
$map[35] = 4
$map[36] = 29

foreach ($targetRow in $map.Keys) {
    $sourceRow = $map[$targetRow]
    $ws.Cells.Item($targetRow, 4).Value = $snapD[$sourceRow]
    $ws.Cells.Item($targetRow, 13).Value = $snapM[$sourceRow]
    $ws.Cells.Item($targetRow, 14).Value = $snapN[$sourceRow]
    $ws.Cells.Item($targetRow, 15).Value = $snapO[$sourceRow]
    $ws.Cells.Item($targetRow, 16).Value = $snapP[$sourceRow]
    $ws.Cells.Item($targetRow, 18).Value = $snapR[$sourceRow]
    $ws.Cells.Item($targetRow, 19).Value = $snapS[$sourceRow]
}

